# RD-Release history.xlsx update
# 1. Row 21 (7.0.6 entry): update release date and add comments about the
#    reconnection / audio-cue fix.
# 2. Row 22: new 7.0.7 release entry (Resolve BRH-11 / AVCTP V1.4).
# 3. Move the active selection to K23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: 7.0.6 release ---------------------------------------------
# Release date moves from 2022-04-29 (44680) to 2022-05-13 (44694)
$ws.Range("C21").Value = 44694
$ws.Range("K21").Value = "Improve the reconnection and audio cue related issues."

# --- Row 22: new 7.0.7 release ------------------------------------------
$ws.Range("A22").Value = "V3 EVT Firmware"
$ws.Range("B22").Value = "7.0.7"
$ws.Range("C22").Value = 44699
$ws.Range("D22").Value = "Zound_Hendrix_M_Lite_V3_hwEVT_btswv7.0.7_20220513"
$ws.Range("E22").Value = "7.0.7"
$ws.Range("G22").Value = 0.6
$ws.Range("H22").Value = 3.1
$ws.Range("J22").Value = "N/A"
$ws.Range("K22").Value = "Resolve the BRH-11. Modify AVCTP to V1.4."

# --- Update the active selection to K23 ---------------------------------
$ws.Range("K23").Select()
